$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D22 description text (T021) ---
$ws.Range("D22").Value = "Vacios plta mqta-protocolo de prueba hermeticidad hidrociclon 1-adecuacion valvula vacio bomba-mqta"

# --- Update Fecha inicio (H column) values ---
$ws.Range("H24").Value = 45775
$ws.Range("H29").Value = 45772
$ws.Range("H30").Value = 45775
$ws.Range("H31").Value = 45775
$ws.Range("H32").Value = 45775
$ws.Range("H33").Value = 45775

# --- Remove the obsolete ID_Transferencia value on row 33 ---
$ws.Range("Q33").ClearContents()

# --- Add new task rows 34 and 35 ---
$ws.Range("D34").Value = "poner en servicio el Hidrociclón 1,  haga la prueba de hermeticidad, proceda a realizar la prueba de vacío."
$ws.Range("A34").Value = "T033"
$ws.Range("B34").Value = "Grupo A"
$ws.Range("C34").Value = "Oscar Rubio"
$ws.Range("E34").Value = "Completada"
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = "Alta"
$ws.Range("H34").Value = 45772
$ws.Range("O33").Copy()
$ws.Range("O34").PasteSpecial(-4122)
$ws.Range("O34").Value = 45775.063888888886
$ws.Range("P34").Value = "[{""fecha"":""2025-04-23"",""estado"":""En curso"",""avance"":10},{""fecha"":""2025-04-25"",""estado"":""En curso"",""avance"":40},{""fecha"":""2025-04-27"",""estado"":""En curso"",""avance"":60}]"

$ws.Range("A35").Value = "T034"
$ws.Range("D35").Value = "poner en servicio el Hidrociclón 1, flushig del sistema -limpieza de filtro interno."
$ws.Range("B35").Value = "Grupo A"
$ws.Range("C35").Value = "Oscar Rubio"
$ws.Range("E35").Value = "En curso"
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = "Alta"
$ws.Range("H35").Value = 45775

$excel.CutCopyMode = 0

# --- Update view / selection state to match author's session ---
$ws.Range("J31").Select()
